# Correct the s1cDNASample metadata sheet so that the "rnaSampleNumber"
# column (C) uses the same plain numeric format as the "s1cDNASampleNumber"
# column (F), instead of free-text "3b N" labels, and update the workbook's
# saved file-path metadata + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column C (rnaSampleNumber) currently holds text like "3b 1".."3b 26".
# Replace each with the plain integer matching the row's sample index
# (same value already present in column F).
for ($row = 2; $row -le 27; $row++) {
    $ws.Cells.Item($row, 3).Value = $row - 1
}

# Update the selection to reflect the newly-normalized column C.
$ws.Range("C2:C27").Select()

# Update the absolute path metadata recorded for this workbook to the new
# database location (the folder these metadata templates now live in).
$wb.AbsPath = "/Users/hollybrown/database_files/s1cDNASample/"

# Turn on iterative calculation with a tightened max-change delta (1E-4)
# to match the corrected calc settings used across the other metadata
# workbooks.
$excel.Iteration = $true
$excel.MaxChange = 0.0001
